$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add a new row of data for "Delegates" -> "Yes", continuing the feature table
$ws.Range("C26").Value = "Delegates"
$ws.Range("E26").Value = "Yes"

# Match the formatting used by the other "Yes" cells in the Supported column
# (bold, green font) by copying the formatting from the row above.
$ws.Range("E25").Copy()
$ws.Range("E26").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Extend the conditional formatting range that colors the "Supported" column
$ws.Range("E10:E25").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E10:E26"))

# Select the cell where the user would naturally land next (below the new row)
$ws.Range("C27").Select()
